$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "['Football', 'soccer', 'PSL', 'score', 'Bafana']"
$ws.Range("E2").Value = "['football,soccer,psl,score,bafana']"

$ws.Range("C3").Value = "['Cricket', 'rugby', 'games', 'athletics', 'race', 'marathon', 'Proteas', 'Springboks', 'Olympics', 'athletics']"
$ws.Range("E3").Value = "['cricket,rugby,games,athletics,race,marathon,proteas,springboks,olympics,athletics']"

$ws.Range("C4").Value = "['Tribal court', 'traditional court', 'imbizo']"
$ws.Range("E4").Value = "['tribal', 'court,traditional', 'court,imbizo']"

$ws.Range("C5").Value = "['Memorial', 'funeral', 'burial', 'after tears', 'mourners']"
$ws.Range("E5").Value = "['memorial,funeral,burial,after', 'tears,mourners']"

$ws.Range("C6").Value = "['Church service', 'ZCC', 'Shembe', 'prayer', 'dedication ']"
$ws.Range("E6").Value = "['church', 'service,zcc,shembe,prayer,dedication']"

$ws.Range("C7").Value = "['Campaign', 'rally', 'door to door', 'register', 'registration']"
$ws.Range("E7").Value = "['campaign,rally,door', 'door,register,registration']"

$ws.Range("C8").Value = "['Debate', 'Various candidates']"
$ws.Range("E8").Value = "['debate,various', 'candidate']"

$ws.Range("C9").Value = "['Vote', 'voting', 'voting station', 'polling station']"
$ws.Range("E9").Value = "['vote,voting,voting', 'station,polling', 'station']"

$ws.Range("C10").Value = "['Taxi', 'taxis', 'taxi drivers', 'over routes']"
$ws.Range("E10").Value = "['taxi,taxis,taxi', 'drivers,over', 'route']"

$ws.Range("C11").Value = "['days of activism', 'against women', 'abuse']"
$ws.Range("E11").Value = "['day', 'activism,against', 'women,abuse']"

$ws.Range("C12").Value = "['Career', 'hobby', 'carnival', 'exhibition', 'inauguration', 'unveiling', 'concert', 'DJ']"
$ws.Range("E12").Value = "['career,hobby,carnival,exhibition,inauguration,unveiling,concert,dj']"
